# Workbook currently has data rows 2..199 for the "Rabanito" subset.
# A new weekly record needs to be inserted as row 86 (pushing the former
# rows 86..199 down to 87..200), mirroring how the upstream daily/weekly
# consolidation script prepends newly scraped records to this block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 86 downwards (this also grows the used range / dimension to
# A1:R200 and keeps each cell's formatting, since Excel's row insert carries
# the row's style along).
$ws.Rows("86:86").Insert()

# Populate the newly inserted row 86 with the new record. Most fields are
# identical to the (now shifted) neighbouring rows for this market/product,
# only the date, volume, and weighted average price differ.
$ws.Range("A86").Value = 9
$ws.Range("B86").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C86").Value = "Metropolitana"
$ws.Range("D86").Value = 44546
$ws.Range("E86").Value = 13
$ws.Range("F86").Value = 300000001
$ws.Range("G86").Value = "Rabanito"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 7900
$ws.Range("K86").Value = 2500
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = 2747
$ws.Range("N86").Value = '$/cien unidades (volumen en unidades)'
$ws.Range("O86").Value = "Provincia de Chacabuco"
$ws.Range("P86").Value = 27
$ws.Range("Q86").Value = 100
$ws.Range("R86").Value = "Hortaliza"

# D column holds dates; make sure the inserted row keeps the same date
# number format used throughout column D.
$ws.Range("D86").NumberFormat = $ws.Range("D87").NumberFormat
